$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price / 1h volume change) for the crypto tracker sheet.
# Map of cell address -> new literal text value. These columns store
# plain text (inline strings) in the source data, e.g. "301.30" / "-0.73%",
# so the cell's NumberFormat is forced to Text ("@") before the assignment
# to stop Excel from re-parsing the literal into a float/percentage.
$updates = @{
    "D2" = "301.30";
    "E2" = "-0.73%";
    "D3" = "31.50";
    "E3" = "-1.89%";
    "D4" = "5.092";
    "E4" = "-3.13%";
    "D5" = "0.07367";
    "E5" = "-2.00%";
    "D6" = "2.329";
    "E6" = "54.07%";
    "D7" = "7.959";
    "E7" = "0.82%";
    "E8" = "-0.77%";
    "D9" = "0.9181";
    "E9" = "-0.45%";
    "D10" = "0.1713";
    "E10" = "0.84%";
    "D11" = "0.07601";
    "E11" = "-4.62%";
    "D12" = "0.08059";
    "E12" = "0.47%";
    "D13" = "0.02990";
    "E13" = "-1.19%";
    "D14" = "0.09932";
    "E14" = "0.04%";
    "D15" = "0.001504";
    "E15" = "0.45%";
    "D16" = "0.006094";
    "E16" = "-6.20%";
    "E17" = "-0.17%";
    "D18" = "2.221";
    "E18" = "-0.49%";
    "E19" = "-0.23%";
    "D20" = "0.1310";
    "E20" = "-2.51%";
    "D21" = "4.652";
    "E21" = "3.43%";
    "D22" = "0.04634";
    "E22" = "0.90%";
    "E23" = "-3.26%";
    "E24" = "0.81%";
    "D25" = "0.004484";
    "E25" = "0.83%";
    "E26" = "-7.06%";
    "E27" = "49.91%";
    "D39" = "0.01737";
    "E39" = "1.42%";
    "D40" = "0.04507";
    "E40" = "0.37%";
    "D41" = "0.007198";
    "E41" = "3.31%";
    "E42" = "-0.19%";
    "E43" = "4.30%";
    "E44" = "-16.60%";
    "D45" = "0.00006305";
    "E45" = "2.08%";
    "D46" = "0.8085";
    "E46" = "-56.58%";
    "E47" = "-33.32%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
